$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - force text format to avoid numeric auto-conversion
$dCells = @("D2", "D3", "D4", "D5", "D7", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($cell in $dCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.124.76"
$ws.Range("D3").Value = "1.897.98"
$ws.Range("D4").Value = "1.003"
$ws.Range("D5").Value = "306.93"
$ws.Range("D7").Value = "0.5234"
$ws.Range("D10").Value = "21.35"
$ws.Range("D11").Value = "0.9024"
$ws.Range("D12").Value = "0.08168"
$ws.Range("D13").Value = "95.26"
$ws.Range("D14").Value = "1.852.37"
$ws.Range("D15").Value = "5.351"
$ws.Range("D17").Value = "0.000008646"
$ws.Range("D18").Value = "14.68"
$ws.Range("D20").Value = "27.167.50"
$ws.Range("D21").Value = "5.113"
$ws.Range("D22").Value = "10.79"
$ws.Range("D23").Value = "6.464"
$ws.Range("D24").Value = "2.333"
$ws.Range("D25").Value = "149.09"
$ws.Range("D26").Value = "18.22"
$ws.Range("D27").Value = "1.745"
$ws.Range("D28").Value = "115.54"
$ws.Range("D30").Value = "4.890"
$ws.Range("D31").Value = "0.09212"
$ws.Range("D33").Value = "0.7927"
$ws.Range("D34").Value = "1.219"
$ws.Range("D35").Value = "2.981"
$ws.Range("D36").Value = "3.362"
$ws.Range("D37").Value = "2.650"
$ws.Range("D38").Value = "0.5677"
$ws.Range("D39").Value = "0.01993"
$ws.Range("D41").Value = "9.033"
$ws.Range("D42").Value = "6.587"
$ws.Range("D43").Value = "116.46"
$ws.Range("D44").Value = "0.1511"
$ws.Range("D45").Value = "0.4873"
$ws.Range("D47").Value = "10.17"
$ws.Range("D49").Value = "38.33"
$ws.Range("D50").Value = "63.91"
$ws.Range("D51").Value = "0.05958"

foreach ($cell in $dCells) {
    $ws.Range($cell).Style = "Normal"
}

# Volume(1h) column (E) updates - plain text percentage strings
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  +0.51%  "
